# Common: Added maximum power for a mod
# Adds a new "power" column (C) to the "mods" worksheet with values for
# each existing mod row, mirroring the header style used for the other
# header cells in row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mods")

# New header cell for the "power" column, styled like the existing headers
$ws.Range("C1").Value = "power"
$ws.Range("C1").Style = $ws.Range("A1").Style

# Maximum power values for each existing mod row
$ws.Range("C2").Value = 100
$ws.Range("C3").Value = 25
$ws.Range("C4").Value = 70
$ws.Range("C5").Value = 50
$ws.Range("C6").Value = 80

# Widen the new column to fit its content
$ws.Columns.Item(3).ColumnWidth = 18.5

# Reflect the new active selection on the sheet
[void]$ws.Range("C6").Select()
